$wb = $excel.ActiveWorkbook

# The original "Sheet1" is replaced by a new "ValidLogin" sheet (this is why
# the saved file shows sheetId jump from 1 -> 2, and the worksheet part gets
# a fresh xr:uid). Add the new sheet first, rename it, then drop the old one.
$ws = $wb.Worksheets.Add()
$ws.Name = "ValidLogin"
[void]$wb.Worksheets.Item("Sheet1").Delete()

# Row 1: column headers for the login form.
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Row 2: sample credentials used to validate the login page.
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# View state captured in the saved workbook: zoom level and active selection.
$excel.ActiveWindow.Zoom = 175
[void]$ws.Range("B3").Select()
